$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-11 for columns I (I0) and J (IF)
$values = @{
    2  = @(7, 8)
    3  = @(5, 7)
    4  = @(7, 8)
    5  = @(9, 9)
    6  = @(1, 1)
    7  = @(9, 9)
    8  = @(7, 7)
    9  = @(8, 8)
    10 = @(2, 4)
    11 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
